$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet holds a running daily log in columns A (Date) / B (MarketObjects),
# with the last existing row being row 75 (2025-08-25). Append three more
# days (2025-08-26 .. 2025-08-28), each with an empty market-objects list.

$newRows = @(
    @("2025-08-26", "[]"),
    @("2025-08-27", "[]"),
    @("2025-08-28", "[]")
)

$startRow = 76
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i

    $dateCell = $ws.Cells.Item($r, 1)
    # Force the date-shaped text to be stored as a literal string rather
    # than letting Excel auto-convert it into a date serial value, then
    # strip the temporary text format so the cell keeps the default style
    # (matching the other plain data rows in the sheet).
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $newRows[$i][0]
    $dateCell.ClearFormats()

    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}
